$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 35, shifting existing rows 35-94 down to 36-95
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the latest transaction entry
$ws.Range("R35").Value = "balance your axis"
$ws.Range("S35").Value = "2024-09-09 08:10:16"
